# Update countries & provincias Spain
# Refreshes the COVID-19 country data snapshot: updates the "last updated"
# timestamp, refreshes case counts for several countries, and re-ranks a
# few country pairs whose totals crossed over (new leader written into the
# higher row, the other country's figures shifted down unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Junio de 2020 a las 00:19"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 2631110
$ws.Range("C4").Value = 34573
$ws.Range("D4").Value = 1090209
$ws.Range("E4").Value = 1412504
$ws.Range("G4").Value = 245
$ws.Range("H4").Value = 128397

# Row 5: 'Brasil' -> 'Brasil'
$ws.Range("B5").Value = 1344143
$ws.Range("C5").Value = 28202
$ws.Range("D5").Value = 733848
$ws.Range("E5").Value = 552673
$ws.Range("G5").Value = 519
$ws.Range("H5").Value = 57622

# Row 10: 'Peru' -> 'Peru'
$ws.Range("B10").Value = 279419
$ws.Range("C10").Value = 3430
$ws.Range("D10").Value = 167998
$ws.Range("E10").Value = 102104
$ws.Range("G10").Value = 182
$ws.Range("H10").Value = 9317

# Row 24: 'Colombia' -> 'Colombia'
$ws.Range("B24").Value = 91769
$ws.Range("C24").Value = 3178
$ws.Range("D24").Value = 38280
$ws.Range("E24").Value = 50383
$ws.Range("G24").Value = 167
$ws.Range("H24").Value = 3106

# Row 31: 'Ecuador' -> 'Ecuador'
$ws.Range("B31").Value = 55255
$ws.Range("C31").Value = 681
$ws.Range("D31").Value = 27058
$ws.Range("E31").Value = 23768
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = 4429

# Row 49: 'Irlanda' -> 'Barein'
$ws.Range("A49").Value = "Barein"
$ws.Range("B49").Value = 25705
$ws.Range("C49").Value = 438
$ws.Range("D49").Value = 20517
$ws.Range("E49").Value = 5105
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 83

# Row 50: 'Barein' -> 'Irlanda'
$ws.Range("A50").Value = "Irlanda"
$ws.Range("B50").Value = 25439
$ws.Range("C50").Value = 2
$ws.Range("D50").Value = 23364
$ws.Range("E50").Value = 340
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 1735

# Row 89: 'Republica de Yibuti' -> 'Bulgaria'
$ws.Range("A89").Value = "Bulgaria"
$ws.Range("B89").Value = 4691
$ws.Range("C89").Value = 66
$ws.Range("D89").Value = 2508
$ws.Range("E89").Value = 1964
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 219

# Row 90: 'Bulgaria' -> 'Republica de Yibuti'
$ws.Range("A90").Value = "Republica de Yibuti"
$ws.Range("B90").Value = 4643
$ws.Range("D90").Value = 4348
$ws.Range("E90").Value = 243
$ws.Range("H90").Value = 52

# Row 133: 'Niger' -> 'Niger'
$ws.Range("B133").Value = 1074
$ws.Range("C133").Value = 12
$ws.Range("D133").Value = 939
$ws.Range("E133").Value = 68

# Row 155: 'Montenegro' -> 'Surinam'
$ws.Range("A155").Value = "Surinam"
$ws.Range("B155").Value = 490
$ws.Range("C155").Value = 57
$ws.Range("D155").Value = 199
$ws.Range("E155").Value = 280
$ws.Range("H155").Value = 11

# Row 156: 'Surinam' -> 'Montenegro'
$ws.Range("A156").Value = "Montenegro"
$ws.Range("B156").Value = 481
$ws.Range("C156").Value = 12
$ws.Range("D156").Value = 315
$ws.Range("E156").Value = 157
$ws.Range("H156").Value = 9

# Row 209: 'Groenlandia' -> 'Islas Malvinas'
$ws.Range("A209").Value = "Islas Malvinas"

# Row 210: 'Islas Malvinas' -> 'Groenlandia'
$ws.Range("A210").Value = "Groenlandia"
